# Weekly update for Naranja (Vega Central Mapocho de Santiago) sheet:
# a new price record is inserted as row 1162, pushing the existing
# rows 1162-1204 down to 1163-1205.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 1162 (shifts 1162..1204 -> 1163..1205).
$ws.Rows.Item(1162).Insert()

# Populate the freshly inserted row 1162 with the new weekly record.
# (K/L/Q/R/T repeat the same "Fukumoto / Primera / $/caja 18 kilos granel /
# Región de O'Higgins / 18" combination as the record that used to sit here.)
$ws.Range("A1162").Value = 9
$ws.Range("B1162").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C1162").Value = "Metropolitana"
$ws.Range("D1162").Value = 45075
$ws.Range("E1162").Value = 13
$ws.Range("F1162").Value = "Fruta"
$ws.Range("G1162").Value = 100102
$ws.Range("H1162").Value = "Cítricos"
$ws.Range("I1162").Value = 100102005
$ws.Range("J1162").Value = "Naranja"
$ws.Range("K1162").Value = "Fukumoto"
$ws.Range("L1162").Value = "Primera"
$ws.Range("M1162").Value = 470
$ws.Range("N1162").Value = 10500
$ws.Range("O1162").Value = 11000
$ws.Range("P1162").Value = 10766
$ws.Range("Q1162").Value = "$/caja 18 kilos granel"
$ws.Range("R1162").Value = "Región de O'Higgins"
$ws.Range("S1162").Value = 598
$ws.Range("T1162").Value = 18
